# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetExhibition.Range("F3").Value = 1422
$sheetExhibition.Range("F4").Value = 1059
$sheetExhibition.Range("F8").Value = 215
$sheetExhibition.Range("F14").Value = 377
$sheetExhibition.Range("F18").Value = 403
$sheetExhibition.Range("F27").Value = 226

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F4").Value = 1422
$sheetAll.Range("F5").Value = 1059
$sheetAll.Range("F12").Value = 215
$sheetAll.Range("F19").Value = 377
$sheetAll.Range("F23").Value = 403
$sheetAll.Range("F39").Value = 226
